$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. B12 weight changes from 4 to 3
$ws.Cells.Item(12, 2).Value = 3

# 2. Insert a new row at 38 ("ACI", 2), shifting old rows 38-43 down to 39-44
$ws.Range("A38").EntireRow.Insert()
$ws.Cells.Item(38, 1).Value = "ACI"
$ws.Cells.Item(38, 2).Value = 2

# Copy formatting from the row below (now row 39) onto the new row 38 so the
# inserted cells match the sheet's existing label style.
$ws.Range("A39:B39").Copy()
$ws.Range("A38:B38").PasteSpecial(-4122)

# 3. Update the (shifted) RECKITTBEN row's weight from 3 to 2
$ws.Cells.Item(44, 2).Value = 2

# 4. Append a brand new row 45 for JMISMDL, 1
$ws.Cells.Item(45, 1).Value = "JMISMDL"
$ws.Cells.Item(45, 2).Value = 1

# Match formatting of the new JMISMDL label cell to the other label cells
$ws.Range("A44:B44").Copy()
$ws.Range("A45:B45").PasteSpecial(-4122)
